$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.426.92"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").Value = "1.617.14"
$ws.Range("E3").Value = "  +1.53%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.99"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +0.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.17"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.25%  "
$ws.Range("D12").Value = "1.843.82"
$ws.Range("E12").Value = "  +1.51%  "
$ws.Range("D13").Value = "1.625.43"
$ws.Range("E13").Value = "  +2.02%  "
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("E15").Value = "  +0.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.83"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "236.54"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +10.01%  "
$ws.Range("D18").Value = "26.428.97"
$ws.Range("E18").Value = "  +0.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.79"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +5.46%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.09"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.03%  "
$ws.Range("E24").Value = "  +3.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.08"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.66%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E28").Value = "  +0.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.48"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.39%  "
$ws.Range("E30").Value = "  +0.40%  "
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("D32").Value = "1.528.43"
$ws.Range("E32").Value = "  +7.34%  "
$ws.Range("E33").Value = "  +1.50%  "
$ws.Range("E34").Value = "  +0.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.52"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +3.94%  "
$ws.Range("E37").Value = "  -0.33%  "
$ws.Range("E38").Value = "  +0.40%  "
$ws.Range("E39").Value = "  +0.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.93"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.79%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("E42").Value = "  +1.65%  "
$ws.Range("D43").Value = "1.755.05"
$ws.Range("E43").Value = "  +1.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.760"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.911"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.46"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.54"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +4.50%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0103"
$ws.Range("E48").Value = "  -0.93%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.50"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.00%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0502"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0958"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.70%  "
